# Apply the Sep 5 2023 crypto price/volume refresh (and one rank swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "1.002", "0.5550")
    # are not silently coerced to Number (which would drop trailing zeros).
    $rng.NumberFormat = "@"
    $rng.Value = $text
    # Drop back to the default style so we do not leave stray formatting
    # on cells that originally had none.
    $rng.Style = "Normal"
}

Set-TextCell "D2" "25.766.98"
Set-TextCell "E2" "  -0.96%  "

Set-TextCell "D3" "1.625.46"
Set-TextCell "E3" "  -1.06%  "

Set-TextCell "D4" "1.001"
Set-TextCell "E4" "  -0.50%  "

Set-TextCell "D5" "214.66"
Set-TextCell "E5" "  -0.45%  "

Set-TextCell "D6" "0.5071"
Set-TextCell "E6" "  -1.20%  "

Set-TextCell "E7" "  -0.49%  "

Set-TextCell "D8" "0.2558"
Set-TextCell "E8" "  -1.05%  "

Set-TextCell "D9" "0.06339"
Set-TextCell "E9" "  -0.67%  "

Set-TextCell "D10" "19.35"
Set-TextCell "E10" "  -2.26%  "

Set-TextCell "D11" "0.07762"
Set-TextCell "E11" "  -0.65%  "

Set-TextCell "D12" "4.255"
Set-TextCell "E12" "  -0.85%  "

Set-TextCell "D13" "1.627.74"
Set-TextCell "E13" "  -1.12%  "

Set-TextCell "D14" "1.849.49"
Set-TextCell "E14" "  -1.15%  "

Set-TextCell "D15" "0.5550"
Set-TextCell "E15" "  +1.36%  "

Set-TextCell "D16" "63.60"
Set-TextCell "E16" "  -1.45%  "

Set-TextCell "D17" "0.0₅7477"
Set-TextCell "E17" "  -3.56%  "

Set-TextCell "D18" "25.790.01"
Set-TextCell "E18" "  -1.13%  "

Set-TextCell "D19" "1.002"
Set-TextCell "E19" "  -0.45%  "

Set-TextCell "D20" "4.405"
Set-TextCell "E20" "  -1.08%  "

Set-TextCell "D21" "193.99"
Set-TextCell "E21" "  -2.11%  "

Set-TextCell "D22" "9.758"
Set-TextCell "E22" "  -2.21%  "

Set-TextCell "D23" "5.972"
Set-TextCell "E23" "  -1.84%  "

Set-TextCell "E24" "  -0.69%  "

Set-TextCell "D25" "1.874"
Set-TextCell "E25" "  -1.46%  "

Set-TextCell "E26" "  -1.16%  "

Set-TextCell "D27" "0.1238"
Set-TextCell "E27" "  +1.13%  "

Set-TextCell "D28" "6.733"
Set-TextCell "E28" "  -2.12%  "

Set-TextCell "D29" "15.42"
Set-TextCell "E29" "  -1.71%  "

Set-TextCell "D30" "1.236"
Set-TextCell "E30" "  -0.47%  "

Set-TextCell "E31" "  -0.22%  "

Set-TextCell "D32" "3.298"
Set-TextCell "E32" "  +0.46%  "

Set-TextCell "D33" "3.182"
Set-TextCell "E33" "  -0.80%  "

Set-TextCell "D34" "1.545"
Set-TextCell "E34" "  +0.19%  "

Set-TextCell "E35" "  -1.11%  "

Set-TextCell "D36" "0.8938"
Set-TextCell "E36" "  -2.27%  "

Set-TextCell "B37" "Maker"
Set-TextCell "C37" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D37" "1.132.33"
Set-TextCell "E37" "  +1.51%  "

Set-TextCell "D38" "2.535"
Set-TextCell "E38" "  -2.34%  "

Set-TextCell "B39" "ImmutableX"
Set-TextCell "C39" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D39" "0.5492"
Set-TextCell "E39" "  -0.66%  "

Set-TextCell "D40" "0.01552"
Set-TextCell "E40" "  -1.27%  "

Set-TextCell "D41" "1.000"
Set-TextCell "E41" "  -0.65%  "

Set-TextCell "D42" "5.567"
Set-TextCell "E42" "  +0.82%  "

Set-TextCell "D43" "0.7935"
Set-TextCell "E43" "  -2.23%  "

Set-TextCell "D44" "97.38"
Set-TextCell "E44" "  -2.26%  "

Set-TextCell "D45" "1.773.29"
Set-TextCell "E45" "  -0.59%  "

Set-TextCell "E46" "  -5.97%  "

Set-TextCell "D47" "0.4416"
Set-TextCell "E47" "  -2.89%  "

Set-TextCell "D48" "54.70"
Set-TextCell "E48" "  -1.07%  "

Set-TextCell "D49" "0.05120"
Set-TextCell "E49" "  -3.30%  "

Set-TextCell "D50" "7.566"
Set-TextCell "E50" "  +2.57%  "

Set-TextCell "D51" "1.003"
Set-TextCell "E51" "  -0.38%  "
